$d = $word.ActiveDocument

$find = "Every category is divide in sub category’s four our purpose we need to get an open category A3 certificate (shown in picture X). This allows us to fly or drone autonomously if we can manually control it if needed. If the drone is going to be used by the first aid workers it needs an other certificate but that is not needed for our testing."
$replace = "Every category is divide in sub category’s for our purpose we need to get an open category A3 certificate (shown in picture X). This allows us to fly or drone autonomously if we can manually control it if needed. If the drone is going to be used by the first aid workers it needs an other certificate but that is not needed for our testing."

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
